# Fruta / hortaliza, semanal
#
# The published dataset gained one new weekly observation. A new row is
# inserted at row 50 (pushing every following row down by one), and the
# sheet's last existing row is preserved at the very end as the new row 151
# (Excel's native row-insert behaviour already shifts rows 50-150 down into
# 51-151, so only the brand-new row's values need to be written).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 50; everything from the old row 50
# through the old row 150 shifts down by one (old row 150 -> new row 151).
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new weekly observation.
$ws.Cells.Item(50, 1).Value = 1
$ws.Cells.Item(50, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(50, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(50, 4).Value = 45002
$ws.Cells.Item(50, 5).Value = 15
$ws.Cells.Item(50, 6).Value = "Fruta"
$ws.Cells.Item(50, 7).Value = 100102
$ws.Cells.Item(50, 8).Value = "Cítricos"
$ws.Cells.Item(50, 9).Value = 100102004
$ws.Cells.Item(50, 10).Value = "Mandarina"
$ws.Cells.Item(50, 11).Value = "Murcott"
$ws.Cells.Item(50, 12).Value = "Segunda"
$ws.Cells.Item(50, 13).Value = 300
$ws.Cells.Item(50, 14).Value = 19000
$ws.Cells.Item(50, 15).Value = 20000
$ws.Cells.Item(50, 16).Value = 19500
$ws.Cells.Item(50, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(50, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(50, 19).Value = 975
$ws.Cells.Item(50, 20).Value = 20
